$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old fixture had a 4th throwaway row - drop it, we only keep 3 rows now.
$ws.Rows("4").Delete()

# Refresh the login test data in the first two columns.
$ws.Range("A1").Value = "mngr533455"
$ws.Range("B1").Value = "1q2w3e4r"
$ws.Range("A2").Value = "mngr533450"
$ws.Range("B2").Value = "arunAja"
$ws.Range("A3").Value = "mngr533455"
$ws.Range("B3").Value = "1q2w3e4r"

# The two alternating row-highlight colors are retired - every row now shares
# the same plain Courier New look (no fill), so reset formatting once for the
# whole block instead of keeping per-row colors.
$dataRng = $ws.Range("A1:B3")
$dataRng.ClearFormats()
$dataRng.Font.Name = "Courier New"
$dataRng.VerticalAlignment = -4108

# New third column with a single-letter flag per row.
$ws.Range("C1").Value = "i"
$ws.Range("C2").Value = "v"
$ws.Range("C3").Value = "i"

$ws.Range("C3").Select()
